$d = $word.ActiveDocument

$replacements = @(
    @("2025-01-01 Wednesday", "2025-01-02 Thursday"),
    @("67×65=4355", "16×59=944"),
    @("76×60=4560", "79×56=4424"),
    @("84×70=5880", "68×60=4080"),
    @("18×44=792", "71×35=2485"),
    @("71×43=3053", "17×66=1122"),
    @("80×91=7280", "42×88=3696"),
    @("60×97=5820", "60×50=3000"),
    @("82×82=6724", "98×67=6566"),
    @("79×82=6478", "79×61=4819"),
    @("87×43=3741", "68×77=5236"),
    @("23×32=736", "97×41=3977"),
    @("83×87=7221", "70×95=6650"),
    @("69×58=4002", "82×75=6150"),
    @("83×76=6308", "23×35=805"),
    @("11×20=220", "35×77=2695"),
    @("21×14=294", "94×25=2350"),
    @("11×52=572", "44×30=1320"),
    @("31×22=682", "79×18=1422"),
    @("51×46=2346", "81×88=7128"),
    @("57×56=3192", "32×92=2944"),
    @("64×29=1856", "54×12=648"),
    @("41×52=2132", "53×66=3498"),
    @("97×15=1455", "31×29=899"),
    @("47×98=4606", "74×89=6586"),
    @("93×11=1023", "73×92=6716")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
